$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 5.75
$ws.Range("N5").Value = 4.75
$ws.Range("AC5").Value = 4.75
$ws.Range("AD5").Value = 7
$ws.Range("AE5").Value = 29
$ws.Range("AG5").Value = 9
$ws.Range("AH5").Value = 29
$ws.Range("AI5").Value = 21
$ws.Range("AJ5").Value = 81

# Row 8
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 9

# Row 14
$ws.Range("G14").Value = 1.75
$ws.Range("I14").Value = 4.33
$ws.Range("J14").Value = 2.3
$ws.Range("L14").Value = 4.5
$ws.Range("O14").Value = 1.2
$ws.Range("P14").Value = 4.33
$ws.Range("Q14").Value = 1.65
$ws.Range("R14").Value = 2.2
$ws.Range("U14").Value = 1.62
$ws.Range("V14").Value = 2.2
$ws.Range("W14").Value = 9
$ws.Range("X14").Value = 9.5
$ws.Range("AH14").Value = 23
$ws.Range("AL14").Value = 34
$ws.Range("AO14").Value = 9
$ws.Range("AQ14").Value = 26
$ws.Range("AW14").Value = 6.5
$ws.Range("AY14").Value = 26
$ws.Range("BA14").Value = 81

# Row 15
$ws.Range("G15").Value = 3.25
$ws.Range("H15").Value = 2.9
$ws.Range("I15").Value = 2.4
$ws.Range("L15").Value = 3.25
$ws.Range("M15").Value = 1.13
$ws.Range("N15").Value = 6
$ws.Range("O15").Value = 1.53
$ws.Range("P15").Value = 2.38
$ws.Range("Q15").Value = 2.7
$ws.Range("R15").Value = 1.44
$ws.Range("S15").Value = 1.62
$ws.Range("T15").Value = 2.2
$ws.Range("U15").Value = 2.2
$ws.Range("V15").Value = 1.62
$ws.Range("W15").Value = 7.5
$ws.Range("AC15").Value = 6
$ws.Range("AF15").Value = 81
$ws.Range("AI15").Value = 11
$ws.Range("AJ15").Value = 23
$ws.Range("AT15").Value = 2.2
$ws.Range("AU15").Value = 9.5

# Row 17
$ws.Range("I17").Value = 3
$ws.Range("M17").Value = 1.17
$ws.Range("N17").Value = 5
$ws.Range("AE17").Value = 21
$ws.Range("AG17").Value = 6.5
$ws.Range("AO17").Value = 17
$ws.Range("AP17").Value = 34
$ws.Range("AU17").Value = 10
$ws.Range("AZ17").Value = 67
$ws.Range("BA17").Value = 126

# Row 18
$ws.Range("G18").Value = 1.42
$ws.Range("I18").Value = 7.5
$ws.Range("K18").Value = 2.1
$ws.Range("L18").Value = 8.5
$ws.Range("AH18").Value = 34
$ws.Range("AI18").Value = 23
$ws.Range("AL18").Value = 67
$ws.Range("AP18").Value = 23
$ws.Range("AQ18").Value = 23
$ws.Range("AW18").Value = 9
$ws.Range("AX18").Value = 41
$ws.Range("AZ18").Value = 251
$ws.Range("BA18").Value = 301

# Row 19
$ws.Range("G19").Value = 2.95
$ws.Range("H19").Value = 3.25
$ws.Range("I19").Value = 2.25
$ws.Range("J19").Value = 3.5
$ws.Range("K19").Value = 2.07
$ws.Range("L19").Value = 2.82
$ws.Range("N19").Value = 9.8
$ws.Range("Q19").Value = 1.87
$ws.Range("W19").Value = 9.5
$ws.Range("X19").Value = 15.5
$ws.Range("Y19").Value = 10.5
$ws.Range("Z19").Value = 37
$ws.Range("AA19").Value = 25
$ws.Range("AC19").Value = 10
$ws.Range("AD19").Value = 6.3
$ws.Range("AE19").Value = 13
$ws.Range("AG19").Value = 8.25
$ws.Range("AH19").Value = 11.5
$ws.Range("AJ19").Value = 23
$ws.Range("AK19").Value = 17.5
$ws.Range("AN19").Value = 4.9
$ws.Range("AO19").Value = 16
$ws.Range("AP19").Value = 23
$ws.Range("AQ19").Value = 75
$ws.Range("AS19").Value = 300
$ws.Range("AT19").Value = 2.55
$ws.Range("AV19").Value = 60
$ws.Range("AW19").Value = 4.2
$ws.Range("AX19").Value = 11.75
$ws.Range("AY19").Value = 19
$ws.Range("AZ19").Value = 45

# Row 20
$ws.Range("G20").Value = 2.55
$ws.Range("I20").Value = 2.65
$ws.Range("J20").Value = 3.15
$ws.Range("L20").Value = 3.25
$ws.Range("P20").Value = 2.47
$ws.Range("Q20").Value = 2.2
$ws.Range("W20").Value = 6.9
$ws.Range("X20").Value = 11.5
$ws.Range("Y20").Value = 10.25
$ws.Range("Z20").Value = 27
$ws.Range("AA20").Value = 25
$ws.Range("AB20").Value = 40
$ws.Range("AD20").Value = 6.1
$ws.Range("AE20").Value = 17
$ws.Range("AG20").Value = 7
$ws.Range("AH20").Value = 11.75
$ws.Range("AI20").Value = 10.5
$ws.Range("AJ20").Value = 29
$ws.Range("AK20").Value = 26
$ws.Range("AL20").Value = 45
$ws.Range("AN20").Value = 4.3
$ws.Range("AO20").Value = 13.5
$ws.Range("AP20").Value = 24
$ws.Range("AQ20").Value = 60
$ws.Range("AR20").Value = 110
$ws.Range("AU20").Value = 7.5
$ws.Range("AV20").Value = 75
$ws.Range("AW20").Value = 4.4
$ws.Range("AX20").Value = 14
$ws.Range("AZ20").Value = 65
